# Refresh the cryptos price/volume snapshot to match the latest scrape.
# Most cells are simple text replacements; a handful of Price cells look
# like plain decimals (single '.'), so Excel's smart-typing would otherwise
# coerce them to numbers.  For those we briefly force a Text number format,
# assign the literal string, then restore the default 'Normal' style so no
# stray cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.290.72'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '3.020.01'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '3.013.50'
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.09'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D15").Value = '3.519.22'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '61.409.13'
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").Value = '3.027.64'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.681'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.91%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.25%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '55.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("B35").Value = 'Stacks'
$ws.Range("C35").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '462.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.92%  '
$ws.Range("D38").Value = '3.231.63'
$ws.Range("E38").Value = '  +4.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0792'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0388'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.118'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.248'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.72%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.109'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").Value = '0.0₃0498'
$ws.Range("E50").Value = '  -7.95%  '
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.46%  '
